# Add data for 2021-09-20: one additional day (September 12) of carjacking
# incidents gets folded into the "Through" workbook. This bumps the
# "current month to date" column (B) for the neighborhoods that had a new
# incident, and also the historical "September" column for the matching
# year(s) since those are also year-to-date-through-this-date comparisons.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet name + page title reflect the new "through" date.
$ws.Name = "Through 2021-09-12"
$ws.Range("B1").Value = "September 2021 (through September 12)"

# Garfield Park: September 2019 (T) 1 -> 2
$ws.Range("T2").Value = 2

# Humboldt Park: September 2018 (AC) new +1, September 2016 (AU) new +1
$ws.Range("AC4").Value = 1
$ws.Range("AU4").Value = 1

# Austin: September 2018 (AC) 2 -> 3
$ws.Range("AC5").Value = 3

# Roseland: September 2021 through-date (B) 3 -> 4
$ws.Range("B6").Value = 4

# Auburn Gresham: September 2021 through-date (B) 4 -> 5
$ws.Range("B7").Value = 5

# Little Village: September 2017 (AL) new +1
$ws.Range("AL9").Value = 1

# West Town: September 2021 through-date (B) 1 -> 2
$ws.Range("B10").Value = 2

# South Chicago: September 2020 (K) new +1
$ws.Range("K22").Value = 1

# Logan Square: September 2018 (AC) new +1
$ws.Range("AC28").Value = 1

# West Loop: September 2020 (K) 2 -> 3
$ws.Range("K31").Value = 3

# Chicago Lawn: September 2019 (T) new +1
$ws.Range("T32").Value = 1

# Hyde Park: September 2017 (AL) new +1
$ws.Range("AL34").Value = 1

# West Elsdon: September 2021 through-date (B) new +1
$ws.Range("B40").Value = 1

# Washington Heights: September 2021 through-date (B) new +1
$ws.Range("B41").Value = 1

# Woodlawn: September 2019 (T) new +1
$ws.Range("T43").Value = 1

# Grand Crossing: September 2018 (AC) 2 -> 3
$ws.Range("AC55").Value = 3

# Clearing: September 2017 (AL) new +1
$ws.Range("AL66").Value = 1

# Near South Side: September 2015 (BD) new +1
$ws.Range("BD85").Value = 1

# Uptown: September 2016 (AU) 1 -> 2
$ws.Range("AU97").Value = 2
